$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2"  = 0.07967938146308559
    "C2"  = 0.8019381115720245
    "D2"  = 1.602517945319166
    "E2"  = 1.265905978072292
    "F2"  = 1.275721982108824

    "B3"  = 0.1996680100171121
    "C3"  = 1.350499864034121
    "D3"  = 3.617881843527005
    "E3"  = 1.902073038431228
    "F3"  = 1.910386039492769

    "B4"  = 0.6016847796089962
    "C4"  = 1.493953804328142
    "D4"  = 4.491473340817588
    "E4"  = 2.119309637787171
    "F4"  = 2.052735548343508

    "B5"  = 0.3601685393389483
    "C5"  = 1.499852406935977
    "D5"  = 4.026777908913376
    "E5"  = 2.006683310568306
    "F5"  = 1.994553887448934
    "G5"  = 49

    "B6"  = 0.554109530473832
    "C6"  = 1.40221327316157
    "D6"  = 3.723570466260433
    "E6"  = 1.929655530466625
    "F6"  = 1.867946810805509
    "G6"  = 48

    "B7"  = 0.3968503597739336
    "C7"  = 1.194979168351206
    "D7"  = 2.489467468001566
    "E7"  = 1.5778046355622
    "F7"  = 1.547043982108814
    "G7"  = 39

    "B8"  = 0.5649144946704114
    "C8"  = 1.22464266767186
    "D8"  = 2.755048971642472
    "E8"  = 1.659834019305085
    "F8"  = 1.581694116082434
    "G8"  = 38

    "B9"  = 0.884415467835616
    "C9"  = 1.104993751378739
    "D9"  = 2.323198242346433
    "E9"  = 1.524204134079957
    "F9"  = 1.27202904791106
    "G9"  = 21

    "B10" = 0.5881942522924533
    "C10" = 0.8478029123914703
    "D10" = 0.9561329667000577
    "E10" = 0.9778205186536318
    "F10" = 0.810614526421041
    "G10" = 14

    "B11" = 0.7033621104908662
    "C11" = 0.8247366046038167
    "D11" = 1.023412291882237
    "E11" = 1.011638419536465
    "F11" = 0.812937600163807
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
